# Apply updated crypto price/volume figures to Sheet1 (columns D and E, rows 2-51)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '36.551.06'
$ws.Range("E2").Value = '  -0.50%  '

$ws.Range("D3").Value = '1.969.69'
$ws.Range("E3").Value = '  +0.65%  '

$ws.Range("E4").Value = '  -0.09%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '243.86'
$ws.Range("E5").Value = '  +0.01%  '

$ws.Range("E6").Value = '  +2.11%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '59.91'
$ws.Range("E7").Value = '  +2.67%  '

$ws.Range("E8").Value = '  -0.01%  '

$ws.Range("E9").Value = '  +1.89%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0789'
$ws.Range("E10").Value = '  -2.25%  '

$ws.Range("E11").Value = '  +0.68%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '14.24'
$ws.Range("E12").Value = '  +4.41%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.843'
$ws.Range("E13").Value = '  +2.90%  '

$ws.Range("D14").Value = '2.263.40'
$ws.Range("E14").Value = '  +0.80%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '21.63'
$ws.Range("E15").Value = '  -2.47%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '5.29'
$ws.Range("E16").Value = '  +0.40%  '

$ws.Range("D17").Value = '1.975.74'
$ws.Range("E17").Value = '  +0.75%  '

$ws.Range("D18").Value = '36.531.21'
$ws.Range("E18").Value = '  -0.48%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '69.89'
$ws.Range("E19").Value = '  +0.30%  '

$ws.Range("D20").Value = '0.0₃0853'
$ws.Range("E20").Value = '  -0.57%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '229.59'
$ws.Range("E21").Value = '  +0.63%  '

$ws.Range("E22").Value = '  -0.69%  '

$ws.Range("E23").Value = '  -0.01%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.44'
$ws.Range("E24").Value = '  +1.51%  '

$ws.Range("E25").Value = '  +1.42%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.144'
$ws.Range("E26").Value = '  +4.49%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.14'

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '162.39'
$ws.Range("E28").Value = '  +1.00%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '19.35'
$ws.Range("E29").Value = '  -0.03%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.33'
$ws.Range("E30").Value = '  +20.03%  '

$ws.Range("E31").Value = '  +1.93%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.80'
$ws.Range("E32").Value = '  +2.81%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0614'
$ws.Range("E33").Value = '  -0.84%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.53'
$ws.Range("E34").Value = '  +7.24%  '

$ws.Range("E35").Value = '  +3.72%  '

$ws.Range("E36").Value = '  -0.01%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.36'
$ws.Range("E37").Value = '  -2.25%  '

$ws.Range("E38").Value = '  +0.42%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.40'
$ws.Range("E39").Value = '  -13.90%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0966'
$ws.Range("E40").Value = '  -3.12%  '

$ws.Range("E41").Value = '  +0.27%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.16'
$ws.Range("E42").Value = '  +0.48%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0209'
$ws.Range("E43").Value = '  -1.27%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '15.89'
$ws.Range("E44").Value = '  -1.14%  '

$ws.Range("D45").Value = '1.366.67'
$ws.Range("E45").Value = '  +1.57%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '88.89'
$ws.Range("E46").Value = '  +1.63%  '

$ws.Range("E47").Value = '  -1.00%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.20'
$ws.Range("E48").Value = '  +0.69%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.82'
$ws.Range("E49").Value = '  -0.51%  '

$ws.Range("E50").Value = '  +5.94%  '

$ws.Range("D51").Value = '2.157.18'
$ws.Range("E51").Value = '  +0.92%  '
